# "Generate Report for Handback" — mark the zh-cn / de-de handback rows as
# complete: update status text, stamp handback datetimes, and fill in the
# "Latest Target File" / "Latest Handback File" columns (with a hyperlink on
# the target-file cell, mirroring the existing source-file hyperlink).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$mdFileName  = "3e1c8fde-3fff-4648-b733-40f34f2333e5.md"
$mdUrl       = "https://github.com/OpenLocalizationTestOrg/oltest/blob/d59fb4d6b42ac23a707513a494ed0f7f37c4fb4c/e2e/3e1c8fde-3fff-4648-b733-40f34f2333e5.md"
$zhXlfName   = "3e1c8fde-3fff-4648-b733-40f34f2333e5.27c636a71e8c1d36de699609486adc26f8d725bb.zh-cn.xlf"
$deXlfName   = "3e1c8fde-3fff-4648-b733-40f34f2333e5.27c636a71e8c1d36de699609486adc26f8d725bb.de-de.xlf"
$statusText  = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: status column mirrors per-language status, so both the
# zh-cn and de-de status cells flip to the new "handed back" wording too.
# ---------------------------------------------------------------------------
$ws1.Range("E2").Value2 = $statusText
$ws1.Range("F2").Value2 = $statusText

# ---------------------------------------------------------------------------
# zh-cn sheet (row 2)
# ---------------------------------------------------------------------------
$ws2.Range("C2").Value2 = $statusText

$ws2Target = $ws2.Range("I2")
$ws2.Hyperlinks.Add($ws2Target, $mdUrl, "", "", $mdFileName) | Out-Null

$ws2.Range("J2").Value2 = $zhXlfName
$ws2.Range("K2").Value2 = "2016-08-12 13:10:37"

# ---------------------------------------------------------------------------
# de-de sheet (row 2)
# ---------------------------------------------------------------------------
$ws3.Range("C2").Value2 = $statusText

$ws3Target = $ws3.Range("I2")
$ws3.Hyperlinks.Add($ws3Target, $mdUrl, "", "", $mdFileName) | Out-Null

$ws3.Range("J2").Value2 = $deXlfName
$ws3.Range("K2").Value2 = "2016-08-12 13:10:47"

# ---------------------------------------------------------------------------
# Column widths: the handoff/handback columns grow wider to comfortably fit
# file names; the engine quantizes ColumnWidth to its character-width grid,
# so the inputs below are chosen to land on the closest achievable width.
# ---------------------------------------------------------------------------
$ws1.Columns.Item(5).ColumnWidth = 29.15   # -> ~29.98 target
$ws1.Columns.Item(6).ColumnWidth = 29.15   # -> ~29.98 target

$ws2.Columns.Item(3).ColumnWidth = 29.15    # -> ~29.98 target
$ws2.Columns.Item(9).ColumnWidth = 38.85    # -> ~39.69 target
$ws2.Columns.Item(10).ColumnWidth = 39.15   # -> 40 target

$ws3.Columns.Item(3).ColumnWidth = 29.15    # -> ~29.98 target
$ws3.Columns.Item(9).ColumnWidth = 38.85    # -> ~39.69 target
$ws3.Columns.Item(10).ColumnWidth = 39.15   # -> 40 target
